# edit.ps1
# Commit: "Refined metadata to be additional tab"
# 1) Adds a new "metadata" worksheet (after "data") describing the panelapp query.
# 2) Refreshes the F-column ("time_taken") timestamps on the "data" sheet to the
#    values captured by the same (re-run) query that produced the metadata sheet.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1. Update the per-row query timestamps on the "data" sheet (column F).
# ---------------------------------------------------------------------------
$newTimestamps = @(
    "2021-10-05 14:33:15.213932",
    "2021-10-05 14:33:15.213940",
    "2021-10-05 14:33:15.213943",
    "2021-10-05 14:33:15.213946",
    "2021-10-05 14:33:15.213949",
    "2021-10-05 14:33:15.213951",
    "2021-10-05 14:33:15.213954",
    "2021-10-05 14:33:15.213957",
    "2021-10-05 14:33:15.213959",
    "2021-10-05 14:33:15.213962",
    "2021-10-05 14:33:15.213965",
    "2021-10-05 14:33:15.213967",
    "2021-10-05 14:33:15.213970",
    "2021-10-05 14:33:15.213972",
    "2021-10-05 14:33:15.213975",
    "2021-10-05 14:33:15.213978",
    "2021-10-05 14:33:15.213980",
    "2021-10-05 14:33:15.213983",
    "2021-10-05 14:33:15.213986",
    "2021-10-05 14:33:15.213988",
    "2021-10-05 14:33:15.213991",
    "2021-10-05 14:33:15.213994",
    "2021-10-05 14:33:15.213996",
    "2021-10-05 14:33:15.213999",
    "2021-10-05 14:33:15.214002",
    "2021-10-05 14:33:15.214005",
    "2021-10-05 14:33:15.214007",
    "2021-10-05 14:33:15.214010",
    "2021-10-05 14:33:15.214012",
    "2021-10-05 14:33:15.214015",
    "2021-10-05 14:33:15.214017",
    "2021-10-05 14:33:15.214020",
    "2021-10-05 14:33:15.214023",
    "2021-10-05 14:33:15.214026",
    "2021-10-05 14:33:15.214028",
    "2021-10-05 14:33:15.214031",
    "2021-10-05 14:33:15.214033",
    "2021-10-05 14:33:15.214036",
    "2021-10-05 14:33:15.214038",
    "2021-10-05 14:33:15.214041",
    "2021-10-05 14:33:15.214044",
    "2021-10-05 14:33:15.214047",
    "2021-10-05 14:33:15.214049",
    "2021-10-05 14:33:15.214052",
    "2021-10-05 14:33:15.214054",
    "2021-10-05 14:33:15.214057",
    "2021-10-05 14:33:15.214059",
    "2021-10-05 14:33:15.214062",
    "2021-10-05 14:33:15.214065",
    "2021-10-05 14:33:15.214067",
    "2021-10-05 14:33:15.214070",
    "2021-10-05 14:33:15.214072",
    "2021-10-05 14:33:15.214075",
    "2021-10-05 14:33:15.214078",
    "2021-10-05 14:33:15.214080",
    "2021-10-05 14:33:15.214083",
    "2021-10-05 14:33:15.214086",
    "2021-10-05 14:33:15.214088",
    "2021-10-05 14:33:15.214091",
    "2021-10-05 14:33:15.214093",
    "2021-10-05 14:33:15.214096",
    "2021-10-05 14:33:15.214098",
    "2021-10-05 14:33:15.214101",
    "2021-10-05 14:33:15.214103",
    "2021-10-05 14:33:15.214107",
    "2021-10-05 14:33:15.214110",
    "2021-10-05 14:33:15.214113",
    "2021-10-05 14:33:15.214116",
    "2021-10-05 14:33:15.214118",
    "2021-10-05 14:33:15.214121",
    "2021-10-05 14:33:15.214123",
    "2021-10-05 14:33:15.214126",
    "2021-10-05 14:33:15.214128",
    "2021-10-05 14:33:15.214131",
    "2021-10-05 14:33:15.214133",
    "2021-10-05 14:33:15.214136",
    "2021-10-05 14:33:15.214141",
    "2021-10-05 14:33:15.214144",
    "2021-10-05 14:33:15.214146",
    "2021-10-05 14:33:15.214149",
    "2021-10-05 14:33:15.214152",
    "2021-10-05 14:33:15.214154",
    "2021-10-05 14:33:15.214157",
    "2021-10-05 14:33:15.214159",
    "2021-10-05 14:33:15.214162",
    "2021-10-05 14:33:15.214165",
    "2021-10-05 14:33:15.214167",
    "2021-10-05 14:33:15.214170",
    "2021-10-05 14:33:15.214172",
    "2021-10-05 14:33:15.214175",
    "2021-10-05 14:33:15.214177",
    "2021-10-05 14:33:15.214180",
    "2021-10-05 14:33:15.214184",
    "2021-10-05 14:33:15.214187",
    "2021-10-05 14:33:15.214190",
    "2021-10-05 14:33:15.214192",
    "2021-10-05 14:33:15.214195",
    "2021-10-05 14:33:15.214198",
    "2021-10-05 14:33:15.214200",
    "2021-10-05 14:33:15.214203",
    "2021-10-05 14:33:15.214205",
    "2021-10-05 14:33:15.214208",
    "2021-10-05 14:33:15.214210",
    "2021-10-05 14:33:15.214213",
    "2021-10-05 14:33:15.214215",
    "2021-10-05 14:33:15.214218",
    "2021-10-05 14:33:15.214220",
    "2021-10-05 14:33:15.214223",
    "2021-10-05 14:33:15.214227",
    "2021-10-05 14:33:15.214230",
    "2021-10-05 14:33:15.214233",
    "2021-10-05 14:33:15.214235",
    "2021-10-05 14:33:15.214238",
    "2021-10-05 14:33:15.214241",
    "2021-10-05 14:33:15.214243",
    "2021-10-05 14:33:15.214245",
    "2021-10-05 14:33:15.214248",
    "2021-10-05 14:33:15.214250",
    "2021-10-05 14:33:15.214253",
    "2021-10-05 14:33:15.214256",
    "2021-10-05 14:33:15.214258",
    "2021-10-05 14:33:15.214260",
    "2021-10-05 14:33:15.214263",
    "2021-10-05 14:33:15.214265",
    "2021-10-05 14:33:15.214268",
    "2021-10-05 14:33:15.214270",
    "2021-10-05 14:33:15.214273",
    "2021-10-05 14:33:15.214275",
    "2021-10-05 14:33:15.214280",
    "2021-10-05 14:33:15.214283",
    "2021-10-05 14:33:15.214286",
    "2021-10-05 14:33:15.214288",
    "2021-10-05 14:33:15.214291",
    "2021-10-05 14:33:15.214293",
    "2021-10-05 14:33:15.214296",
    "2021-10-05 14:33:15.214298",
    "2021-10-05 14:33:15.214301",
    "2021-10-05 14:33:15.214303",
    "2021-10-05 14:33:15.214306",
    "2021-10-05 14:33:15.214308",
    "2021-10-05 14:33:15.214311",
    "2021-10-05 14:33:15.214313",
    "2021-10-05 14:33:15.214316",
    "2021-10-05 14:33:15.214318",
    "2021-10-05 14:33:15.214321",
    "2021-10-05 14:33:15.214323"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# ---------------------------------------------------------------------------
# 2. Add the new "metadata" worksheet, placed right after "data".
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row (B1:G1) - reuse the bold/bordered/centered header style already
# used by the "data" sheet's header row.
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$dataSheet.Range("B1:F1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Data row 2.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Autism"
$ws.Range("C2").Value = 51
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.168"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "2021-09-21T04:49:10.848940Z"
$ws.Range("F2").Value = "2021-10-05 14:33:15.210321"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/51/?format=json"

$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

Write-Host "metadata sheet added; data sheet timestamps refreshed"
